$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.275.76"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.025.15"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.28"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.03"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.022.23"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  +5.72%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.55"
$ws.Range("E14").Value = "  +5.03%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.261.78"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.526.37"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("E19").Value = "  +18.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.022.73"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.40"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.38"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.97"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.78"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.45"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +6.74%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0998"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.24"
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.48"
$ws.Range("E38").Value = "  +10.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.44"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  -4.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.61"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0359"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "378.96"
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.706.84"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.95"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.50"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  +2.69%  "
